# Add PF/1.0.3 to meta-sheet
# Appends a new row (row 3) below the existing "PF/1.0.0" row with the
# new version label in column A and "X" markers in columns B-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
